# feat: Support blog static field
# Adds a new "Static" column (D) to the "Blog" sheet, flagging which
# blog categories are static (DisplayName/Description rows 2, 4, 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blog")

# New header for column D
$ws.Range("D1").Value = "Static"

# Flag the static categories with a value of 1
$ws.Range("D2").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1

# Match the cursor/selection left behind by the edit
$ws.Range("M16").Select()
